$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the old values for rows 2, 3, 4 (columns D, J, K, L, M, P)
$old2D = $ws.Range("D2").Value()
$old2J = $ws.Range("J2").Value()
$old2K = $ws.Range("K2").Value()
$old2L = $ws.Range("L2").Value()
$old2M = $ws.Range("M2").Value()
$old2P = $ws.Range("P2").Value()

$old3D = $ws.Range("D3").Value()
$old3J = $ws.Range("J3").Value()
$old3K = $ws.Range("K3").Value()
$old3L = $ws.Range("L3").Value()
$old3M = $ws.Range("M3").Value()
$old3P = $ws.Range("P3").Value()

$old4D = $ws.Range("D4").Value()
$old4J = $ws.Range("J4").Value()
$old4K = $ws.Range("K4").Value()
$old4L = $ws.Range("L4").Value()
$old4M = $ws.Range("M4").Value()
$old4P = $ws.Range("P4").Value()

# Row 2 takes the old row 4 values
$ws.Range("D2").Value = $old4D
$ws.Range("J2").Value = $old4J
$ws.Range("K2").Value = $old4K
$ws.Range("L2").Value = $old4L
$ws.Range("M2").Value = $old4M
$ws.Range("P2").Value = $old4P

# Row 3 takes the old row 2 values
$ws.Range("D3").Value = $old2D
$ws.Range("J3").Value = $old2J
$ws.Range("K3").Value = $old2K
$ws.Range("L3").Value = $old2L
$ws.Range("M3").Value = $old2M
$ws.Range("P3").Value = $old2P

# Row 4 takes the old row 3 values
$ws.Range("D4").Value = $old3D
$ws.Range("J4").Value = $old3J
$ws.Range("K4").Value = $old3K
$ws.Range("L4").Value = $old3L
$ws.Range("M4").Value = $old3M
$ws.Range("P4").Value = $old3P
